# [GRE] Realizando revisões 4, 5 e 6 de requisitos
#
# Adds three new worksheets ("Revisão 4", "Revisão 5", "Revisão 6") right
# after "Revisão 3" and before "Legenda", each built as a copy of the
# "Revisão 3" template (so formatting / column widths / page margins match),
# then overwrites column A with the new requirement codes and blanks out the
# trailing (8th) template row. Finally "Revisão 6" is left as the active /
# selected sheet, matching the committed workbook state.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Revisão 3")

# ---------------------------------------------------------------------
# Revisão 4 (copied right after Revisão 3)
# ---------------------------------------------------------------------
$src.Copy([System.Reflection.Missing]::Value, $src) | Out-Null
$s4 = $wb.Worksheets.Item("Revisão 3 (2)")
$s4.Name = "Revisão 4"

$s4.Range("A5").Value = "RFUN7.5"
$s4.Range("A4").Value = "RFUN7.4"
$s4.Range("A3").Value = "RFUN7.3"
$s4.Range("A2").Value = "RFUN7.2"
$s4.Range("A1").Value = "RFUN7.1"
$s4.Range("A6").Value = "RFUN8.1"
$s4.Range("A7").Value = "RFUN8.2"
$s4.Range("A8:K8").ClearContents() | Out-Null
$s4.Range("A7").Select() | Out-Null

# ---------------------------------------------------------------------
# Revisão 5 (copied right after Revisão 4)
# ---------------------------------------------------------------------
$s4.Copy([System.Reflection.Missing]::Value, $s4) | Out-Null
$s5 = $wb.Worksheets.Item("Revisão 4 (2)")
$s5.Name = "Revisão 5"

$s5.Range("A1").Value = "RFUN9.1"
$s5.Range("A2").Value = "RFUN9.2"
$s5.Range("A3").Value = "RFUN9.3"
$s5.Range("A4").Value = "RFUN9.4"
$s5.Range("A5").Value = "RFUN9.5"
$s5.Range("A6").Value = "RFUN9.6"
$s5.Range("A7").Value = "RFUN9.7"
$s5.Range("A8:K8").ClearContents() | Out-Null
$s5.Range("H10").Select() | Out-Null

# ---------------------------------------------------------------------
# Revisão 6 (copied right after Revisão 5)
# ---------------------------------------------------------------------
$s5.Copy([System.Reflection.Missing]::Value, $s5) | Out-Null
$s6 = $wb.Worksheets.Item("Revisão 5 (2)")
$s6.Name = "Revisão 6"

$s6.Range("A6").Value = "RFUN11.1"
$s6.Range("A7").Value = "RFUN11.2"
$s6.Range("A1").Value = "RFUN10.1"
$s6.Range("A2").Value = "RFUN10.2"
$s6.Range("A3").Value = "RFUN10.3"
$s6.Range("A4").Value = "RFUN10.4"
$s6.Range("A5").Value = "RFUN10.5"
$s6.Range("A8:K8").ClearContents() | Out-Null
$s6.Range("K7").Select() | Out-Null

# Revisão 6 is the sheet that should end up active/selected.
$s6.Activate() | Out-Null
